# chore: update Sheets via scheduled runner
#
# Refreshes computed market/profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets with
# up-to-date values pulled from the market data source.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3383.6667
$ws.Cells.Item(40, 10).Value = 2200.4
$ws.Cells.Item(40, 12).Value = 2200.4
$ws.Cells.Item(40, 14).Value = -2550.4
$ws.Cells.Item(64, 8).Value = 5403.5835
$ws.Cells.Item(64, 10).Value = 8575.75
$ws.Cells.Item(64, 12).Value = 8575.75
$ws.Cells.Item(64, 14).Value = -9071.75
$ws.Cells.Item(67, 8).Value = 5403.5835
$ws.Cells.Item(67, 10).Value = 8575.75
$ws.Cells.Item(67, 12).Value = 8575.75
$ws.Cells.Item(67, 14).Value = -10291.75
$ws.Cells.Item(106, 8).Value = 18687118
$ws.Cells.Item(106, 9).Value = 22423342
$ws.Cells.Item(106, 11).Value = 22423342
$ws.Cells.Item(106, 13).Value = -22422711
$ws.Cells.Item(113, 8).Value = 101740.3
$ws.Cells.Item(113, 9).Value = 112878.336
$ws.Cells.Item(113, 10).Value = 1498
$ws.Cells.Item(113, 11).Value = 112878.336
$ws.Cells.Item(113, 12).Value = 1498
$ws.Cells.Item(113, 13).Value = -109624.336
$ws.Cells.Item(113, 14).Value = -8006
$ws.Cells.Item(123, 8).Value = 91665.89999999999
$ws.Cells.Item(123, 10).Value = 91665.89999999999
$ws.Cells.Item(123, 12).Value = 91665.89999999999
$ws.Cells.Item(123, 14).Value = -101465.9
$ws.Cells.Item(137, 8).Value = 1281.125
$ws.Cells.Item(137, 9).Value = 769.1579
$ws.Cells.Item(137, 10).Value = 1616.5518
$ws.Cells.Item(137, 11).Value = 2307.4737
$ws.Cells.Item(137, 12).Value = 4849.6554
$ws.Cells.Item(137, 13).Value = 242.5263
$ws.Cells.Item(137, 14).Value = -9949.6554
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2881.7036
$ws.Cells.Item(32, 9).Value = 1849.9833
$ws.Cells.Item(32, 10).Value = 5829.476
$ws.Cells.Item(32, 11).Value = 1849.9833
$ws.Cells.Item(32, 12).Value = 5829.476
$ws.Cells.Item(32, 13).Value = -1562.9833
$ws.Cells.Item(32, 14).Value = -6403.476
$ws.Cells.Item(61, 8).Value = 1211.8
$ws.Cells.Item(61, 9).Value = 986.06384
$ws.Cells.Item(61, 10).Value = 2538
$ws.Cells.Item(61, 11).Value = 986.06384
$ws.Cells.Item(61, 12).Value = 2538
$ws.Cells.Item(61, 13).Value = -774.06384
$ws.Cells.Item(61, 14).Value = -2962
$ws.Cells.Item(74, 8).Value = 693.68085
$ws.Cells.Item(74, 9).Value = 539.0323
$ws.Cells.Item(74, 10).Value = 993.3125
$ws.Cells.Item(74, 11).Value = 539.0323
$ws.Cells.Item(74, 12).Value = 993.3125
$ws.Cells.Item(74, 13).Value = 334.9677
$ws.Cells.Item(74, 14).Value = -2741.3125
$ws.Cells.Item(77, 8).Value = 693.68085
$ws.Cells.Item(77, 9).Value = 539.0323
$ws.Cells.Item(77, 10).Value = 993.3125
$ws.Cells.Item(77, 11).Value = 2695.1615
$ws.Cells.Item(77, 12).Value = 4966.5625
$ws.Cells.Item(77, 13).Value = 1672.8385
$ws.Cells.Item(77, 14).Value = -13702.5625
$ws.Cells.Item(123, 8).Value = 35214
$ws.Cells.Item(123, 10).Value = 35214
$ws.Cells.Item(123, 12).Value = 35214
$ws.Cells.Item(123, 14).Value = -45014
$ws.Cells.Item(132, 8).Value = 2171.1875
$ws.Cells.Item(132, 9).Value = 1999.5385
$ws.Cells.Item(132, 11).Value = 5998.6155
$ws.Cells.Item(132, 13).Value = -3468.6155
$ws.Cells.Item(136, 8).Value = 1211.8
$ws.Cells.Item(136, 9).Value = 986.06384
$ws.Cells.Item(136, 10).Value = 2538
$ws.Cells.Item(136, 11).Value = 2958.19152
$ws.Cells.Item(136, 12).Value = 7614
$ws.Cells.Item(136, 13).Value = -408.1915200000003
$ws.Cells.Item(136, 14).Value = -12714
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1212.0714
$ws.Cells.Item(107, 9).Value = 1003.3333
$ws.Cells.Item(107, 10).Value = 1587.8
$ws.Cells.Item(107, 11).Value = 1003.3333
$ws.Cells.Item(107, 12).Value = 1587.8
$ws.Cells.Item(107, 13).Value = 916.6667
$ws.Cells.Item(107, 14).Value = -5427.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).Value = ""
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2682.3562
$ws.Cells.Item(68, 9).Value = 2840.3728
$ws.Cells.Item(68, 10).Value = 2016.4286
$ws.Cells.Item(68, 11).Value = 8521.118399999999
$ws.Cells.Item(68, 12).Value = 6049.2858
$ws.Cells.Item(68, 13).Value = -7710.118399999999
$ws.Cells.Item(68, 14).Value = -7671.2858
$ws.Cells.Item(71, 8).Value = 2682.3562
$ws.Cells.Item(71, 9).Value = 2840.3728
$ws.Cells.Item(71, 10).Value = 2016.4286
$ws.Cells.Item(71, 11).Value = 25563.3552
$ws.Cells.Item(71, 12).Value = 18147.8574
$ws.Cells.Item(71, 13).Value = -21507.3552
$ws.Cells.Item(71, 14).Value = -26259.8574
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2660.6052
$ws.Cells.Item(132, 9).Value = 2028.8966
$ws.Cells.Item(132, 10).Value = 4696.1113
$ws.Cells.Item(132, 11).Value = 6086.6898
$ws.Cells.Item(132, 12).Value = 14088.3339
$ws.Cells.Item(132, 13).Value = -3556.6898
$ws.Cells.Item(132, 14).Value = -19148.3339
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3640.5293
$ws.Cells.Item(61, 9).Value = 4340.5
$ws.Cells.Item(61, 10).Value = 1960.6
$ws.Cells.Item(61, 11).Value = 4340.5
$ws.Cells.Item(61, 12).Value = 1960.6
$ws.Cells.Item(61, 13).Value = -4138.5
$ws.Cells.Item(61, 14).Value = -2364.6
$ws.Cells.Item(113, 8).Value = 3640.5293
$ws.Cells.Item(113, 9).Value = 4340.5
$ws.Cells.Item(113, 10).Value = 1960.6
$ws.Cells.Item(113, 11).Value = 4340.5
$ws.Cells.Item(113, 12).Value = 1960.6
$ws.Cells.Item(113, 13).Value = -2170.5
$ws.Cells.Item(113, 14).Value = -6300.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 500118
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).Value = ""
$ws.Cells.Item(78, 8).Value = 500118
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).Value = ""
$ws.Cells.Item(82, 8).Value = 50000
$ws.Cells.Item(82, 10).Value = 50000
$ws.Cells.Item(82, 12).Value = 50000
$ws.Cells.Item(82, 14).Value = -50766
$ws.Cells.Item(85, 8).Value = 50000
$ws.Cells.Item(85, 10).Value = 50000
$ws.Cells.Item(85, 12).Value = 50000
$ws.Cells.Item(85, 14).Value = -52652
$ws.Cells.Item(122, 8).Value = 2388.889
$ws.Cells.Item(122, 9).Value = 2366.6667
$ws.Cells.Item(122, 10).Value = 2433.3333
$ws.Cells.Item(122, 11).Value = 7100.000100000001
$ws.Cells.Item(122, 12).Value = 7299.999899999999
$ws.Cells.Item(122, 13).Value = -4650.000100000001
$ws.Cells.Item(122, 14).Value = -12199.9999
$ws.Cells.Item(123, 8).Value = 32737.834
$ws.Cells.Item(123, 10).Value = 32737.834
$ws.Cells.Item(123, 12).Value = 32737.834
$ws.Cells.Item(123, 14).Value = -42537.834
$ws.Cells.Item(126, 8).Value = 250925.25
$ws.Cells.Item(126, 9).Value = 333833.66
$ws.Cells.Item(126, 11).Value = 1001500.98
$ws.Cells.Item(126, 13).Value = -999030.98
$ws.Cells.Item(132, 8).Value = 12501730
$ws.Cells.Item(132, 9).Value = 16130189
$ws.Cells.Item(132, 10).Value = 3705.5557
$ws.Cells.Item(132, 11).Value = 48390567
$ws.Cells.Item(132, 12).Value = 11116.6671
$ws.Cells.Item(132, 13).Value = -48388037
$ws.Cells.Item(132, 14).Value = -16176.6671
